# Regenerate orders with updated distance/sizes:
#   D51 -> D55
#   D64 -> D69
#   D80 -> D86
#   S30 -> S31
# Applied to every text cell value in the used range (the substitutions are
# substring replacements that occur inside Condition / Filename_Left /
# Filename_Right / Distance / Size columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value()
        if ($val -is [string]) {
            if ($val.Contains("D51") -or $val.Contains("D64") -or $val.Contains("D80") -or $val.Contains("S30")) {
                $newVal = $val.Replace("D51", "D55").Replace("D64", "D69").Replace("D80", "D86").Replace("S30", "S31")
                if ($newVal -ne $val) {
                    $cell.Value = $newVal
                }
            }
        }
    }
}
